# Add penalty factor for delay in objective function
# Updates pre-computed values in the "Delay Model" data sheet to reflect
# the new objective-function results after adding a delay-penalty term.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("C1").Value = 3437.9810620581961
$ws.Range("D1").Value = 3437.9810620581961
$ws.Range("E1").Value = 3437.9810620581961
$ws.Range("F1").Value = 3437.9810620581961
$ws.Range("O1").Value = 2.1722486486763244
$ws.Range("P1").Value = 0.0042911264477147317
$ws.Range("Q1").Value = 0.0040177754921099425
$ws.Range("R1").Value = 3.6928380058521566

# Row 2
$ws.Range("C2").Value = 7963.3709603826101
$ws.Range("D2").Value = 16639.936190263325
$ws.Range("E2").Value = 7963.3709603826101
$ws.Range("F2").Value = 7963.3709603826101
$ws.Range("O2").Value = 6.7413431515914519
$ws.Range("P2").Value = 0.008793026108856258
$ws.Range("Q2").Value = 0.0079652932464109524
$ws.Range("R2").Value = 2.0414773808484235

# Row 3
$ws.Range("C3").Value = 15719.502686049842
$ws.Range("D3").Value = 24396.067915930555
$ws.Range("E3").Value = 15719.502686049842
$ws.Range("F3").Value = 15719.502686049842
$ws.Range("O3").Value = 9.4478216830169988
$ws.Range("P3").Value = 0.011140260329096346
$ws.Range("Q3").Value = 0.01151572678216432
$ws.Range("R3").Value = 1.4896000059164411

# Row 4
$ws.Range("C4").Value = 23997.701653556047
$ws.Range("D4").Value = 33173.201411883587
$ws.Range("E4").Value = 24494.893075858083
$ws.Range("F4").Value = 24341.823999796783
$ws.Range("L4").Value = 30.029244321927248
$ws.Range("O4").Value = 11.141239101872866
$ws.Range("P4").Value = 0.015376489030264019
$ws.Range("Q4").Value = 0.015269253535713545
$ws.Range("R4").Value = 0.3153986472398807

# Row 5
$ws.Range("C5").Value = 30653.911403578168
$ws.Range("D5").Value = 39821.624078848217
$ws.Range("E5").Value = 31143.31574282272
$ws.Range("F5").Value = 31020.111972763196
$ws.Range("L5").Value = 35.042217294900226
$ws.Range("O5").Value = 13.07890057733656
$ws.Range("P5").Value = 0.019141677999400961
$ws.Range("Q5").Value = 0.018657838274652002
$ws.Range("R5").Value = 0.2080766816274878

# Row 6
$ws.Range("C6").Value = 36296.511994511246
$ws.Range("D6").Value = 45464.224669781288
$ws.Range("E6").Value = 36785.916333755798
$ws.Range("F6").Value = 36473.054073267122
$ws.Range("O6").Value = 13.023724075392526
$ws.Range("P6").Value = 0.021469854452705023
$ws.Range("Q6").Value = 0.022541924953927137
$ws.Range("R6").Value = 0.13485645661499357

# Row 7
$ws.Range("C7").Value = 43452.204522691391
$ws.Range("D7").Value = 55004.535550839893
$ws.Range("E7").Value = 45118.620898540736
$ws.Range("F7").Value = 44104.757740425368
$ws.Range("O7").Value = 10.622044375014756
$ws.Range("P7").Value = 0.026266153946365754
$ws.Range("Q7").Value = 0.025322936705173785
$ws.Range("R7").Value = 0.87471082683117407

# Row 8
$ws.Range("C8").Value = 50546.618019861751
$ws.Range("D8").Value = 69383.614656866615
$ws.Range("E8").Value = 52215.720319221444
$ws.Range("F8").Value = 51717.348890623864
$ws.Range("O8").Value = 18.145465943912708
$ws.Range("P8").Value = 0.028420534942984566
$ws.Range("Q8").Value = 0.028869103576687743
$ws.Range("R8").Value = 4.4307353994038614

# Row 9
$ws.Range("C9").Value = 53755.358450139349
$ws.Range("D9").Value = 72477.823865821629
$ws.Range("E9").Value = 55749.287937245586
$ws.Range("F9").Value = 54774.9327272722
$ws.Range("I9").Value = 55.047657188626353
$ws.Range("O9").Value = 27.766133324155266
$ws.Range("P9").Value = 0.031925172949234942
$ws.Range("Q9").Value = 0.032136515051226992
$ws.Range("R9").Value = 0.18635708509463603

# Row 10
$ws.Range("C10").Value = 59726.204735787484
$ws.Range("D10").Value = 87477.823865821629
$ws.Range("E10").Value = 62895.591850531062
$ws.Range("F10").Value = 60377.114926591887
$ws.Range("O10").Value = 160.7193877507483
$ws.Range("P10").Value = 0.036029988391828825
$ws.Range("Q10").Value = 0.03563005973045491
$ws.Range("R10").Value = 0.10741697000279039

# Row 11 - zeroed out (data moved up / row removed from the effective set)
$ws.Range("A11").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 0
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = 0
$ws.Range("Q11").Value = 0
$ws.Range("R11").Value = 0
